# AICUM Vaccination survey 1.19.2021 - data cleanup
# - Removes an erroneous duplicate "Yes" value from Dean College's
#   "Named local hospital or healthcare provider" column (G11).
# - Normalizes/cleans up a handful of hospital/clinic name entries in the
#   same column for other institutions.
# - Restores the sheet selection/scroll position to the top of the sheet
#   with H14 selected (matching a fresh open/save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dean College (row 11) incorrectly had "Yes" duplicated into the
# "Named local hospital or healthcare provider" column - clear it.
$ws.Range("G11").ClearContents()

# Suffolk (row 20): tidy up institution name.
$ws.Range("G20").Value2 = "UMass Boston"

# Holy Cross (row 24): fix capitalization.
$ws.Range("G24").Value2 = "UMass Medical School"

# Assumption (row 33): shorten the clinic description.
$ws.Range("G33").Value2 = "WP Clinic"

# Williams (row 25): use a semicolon to separate the two named providers.
$ws.Range("G25").Value2 = "Southwester VT Medical Center; Berkshire Medical Center"

# Reset the view: scroll back to the top of the sheet and select H14
# (instead of being scrolled to row 27 with A40 selected).
$ws.Range("H14").Select()
